# Add team record columns (Wins / Losses / Ties) to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1=Wins, AE1=Losses, AF1=Ties ---
# Copy the formatting of an existing header cell (AC1) so the new
# header cells get the same bold/centered/bordered style (style index 1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows 2-51: same team record for every player row ---
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 72   # column AD - Wins
    $ws.Cells.Item($r, 31).Value = 90   # column AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # column AF - Ties
}
